$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.902.84"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.05%  "

$ws.Range("D3").Value = "'1.901.04"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.68%  "

$ws.Range("D4").Value = "'1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").Value = "'0.7642"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.53%  "

$ws.Range("D6").Value = "'240.54"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.66%  "

$ws.Range("E7").Value = "  -0.15%  "

$ws.Range("D8").Value = "'0.3066"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.71%  "

$ws.Range("E9").Value = "  -2.29%  "

$ws.Range("D10").Value = "'0.06849"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.55%  "

$ws.Range("D11").Value = "'0.07979"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.50%  "

$ws.Range("D12").Value = "'1.901.46"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.14%  "

$ws.Range("D13").Value = "'0.7440"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.55%  "

$ws.Range("D14").Value = "'5.164"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.19%  "

$ws.Range("D15").Value = "'91.14"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.06%  "

$ws.Range("D16").Value = "'29.908.15"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.03%  "

$ws.Range("D17").Value = "'13.97"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.74%  "

$ws.Range("D18").Value = "'5.946"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.72%  "

$ws.Range("D19").Value = "'242.94"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.27%  "

$ws.Range("E20").Value = "  -0.46%  "

$ws.Range("D21").Value = "'1.0000"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.19%  "

$ws.Range("D22").Value = "'1.001"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.06%  "

$ws.Range("D23").Value = "'6.950"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.70%  "

$ws.Range("D24").Value = "'166.76"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.51%  "

$ws.Range("D25").Value = "'9.252"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.35%  "

$ws.Range("D26").Value = "'18.73"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.66%  "

$ws.Range("D27").Value = "'0.1295"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.31%  "

$ws.Range("D28").Value = "'2.039"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.70%  "

$ws.Range("D29").Value = "'1.405"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.68%  "

$ws.Range("D30").Value = "'1.518"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.79%  "

$ws.Range("D31").Value = "'4.263"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.63%  "

$ws.Range("D32").Value = "'4.094"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.71%  "

$ws.Range("D33").Value = "'0.05286"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.36%  "

$ws.Range("D34").Value = "'1.254"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.22%  "

$ws.Range("D35").Value = "'0.7284"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.53%  "

$ws.Range("E36").Value = "  -0.31%  "

$ws.Range("D37").Value = "'0.01927"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.56%  "

$ws.Range("D38").Value = "'2.781"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.50%  "

$ws.Range("D39").Value = "'6.185"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.11%  "

$ws.Range("D40").Value = "'0.4421"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.11%  "

$ws.Range("D41").Value = "'72.08"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.98%  "

$ws.Range("E42").Value = "  -0.20%  "

$ws.Range("D43").Value = "'1.889"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.86%  "

$ws.Range("D44").Value = "'0.8309"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.63%  "

$ws.Range("D45").Value = "'7.641"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.89%  "

$ws.Range("D46").Value = "'100.24"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.46%  "

$ws.Range("D47").Value = "'9.792"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.35%  "

$ws.Range("D48").Value = "'2.057.27"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.43%  "

$ws.Range("D49").Value = "'36.15"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.73%  "

$ws.Range("D50").Value = "'1.482"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.34%  "

$ws.Range("D51").Value = "'0.05944"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.14%  "
